$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the Mana/Attack/Health block (B2:D4) with 1, matching the
# (accidental) mass-edit described in the commit.
$ws.Range("B2:D4").Value = 1

# Update the active cell / selection on the sheet.
$ws.Range("G13").Select()
